$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (Price, Volume(1h)) values refreshed for this run
$updates = @(
    @{ Row = 2; Price = '30.581.38'; Volume = '  +0.34%  ' }
    @{ Row = 3; Price = '2.111.14'; Volume = '  +0.11%  ' }
    @{ Row = 4; Price = '1.014'; Volume = '  +1.14%  ' }
    @{ Row = 5; Price = '346.91'; Volume = '  +4.12%  ' }
    @{ Row = 6; Price = '1.012'; Volume = '  +1.04%  ' }
    @{ Row = 7; Price = '0.5256'; Volume = '  +0.14%  ' }
    @{ Row = 8; Price = '0.4515'; Volume = '  -1.30%  ' }
    @{ Row = 9; Price = '54.12'; Volume = '  +1.04%  ' }
    @{ Row = 10; Price = '0.09015'; Volume = '  +0.38%  ' }
    @{ Row = 11; Price = '1.171'; Volume = '  -0.91%  ' }
    @{ Row = 12; Price = '24.38'; Volume = '  -0.36%  ' }
    @{ Row = 13; Price = '2.107.14'; Volume = '  +0.24%  ' }
    @{ Row = 14; Price = '6.804'; Volume = '  +0.03%  ' }
    @{ Row = 15; Price = '8.068'; Volume = '  +2.92%  ' }
    @{ Row = 16; Price = '99.67'; Volume = '  +3.01%  ' }
    @{ Row = 17; Price = '0.00001178'; Volume = '  +4.08%  ' }
    @{ Row = 18; Price = '1.014'; Volume = '  +1.03%  ' }
    @{ Row = 19; Price = '0.06708'; Volume = '  +1.35%  ' }
    @{ Row = 20; Price = '19.33'; Volume = '  +0.05%  ' }
    @{ Row = 21; Price = '1.011'; Volume = '  +1.03%  ' }
    @{ Row = 22; Price = '6.323'; Volume = '  +0.27%  ' }
    @{ Row = 23; Price = '30.656.72'; Volume = '  +0.35%  ' }
    @{ Row = 24; Price = '12.79'; Volume = '  +3.64%  ' }
    @{ Row = 25; Price = '2.387'; Volume = '  +1.09%  ' }
    @{ Row = 26; Price = '2.360.97'; Volume = '  +0.42%  ' }
    @{ Row = 27; Price = '22.34'; Volume = '  -0.14%  ' }
    @{ Row = 28; Price = '165.20'; Volume = '  +1.28%  ' }
    @{ Row = 29; Price = '2.526'; Volume = '  -2.02%  ' }
    @{ Row = 30; Price = '135.97'; Volume = '  +2.33%  ' }
    @{ Row = 31; Price = '1.194'; Volume = '  -0.28%  ' }
    @{ Row = 32; Price = '0.1073'; Volume = '  -0.11%  ' }
    @{ Row = 33; Price = '1.634'; Volume = '  -4.14%  ' }
    @{ Row = 34; Price = '6.353'; Volume = '  +3.16%  ' }
    @{ Row = 35; Price = '3.995'; Volume = '  +1.61%  ' }
    @{ Row = 36; Price = '5.895'; Volume = '  +6.08%  ' }
    @{ Row = 37; Price = '10.20'; Volume = '  -2.34%  ' }
    @{ Row = 38; Price = '0.02641'; Volume = '  +2.51%  ' }
    @{ Row = 39; Price = '0.06834'; Volume = '  +0.08%  ' }
    @{ Row = 40; Price = '0.2323'; Volume = '  +1.38%  ' }
    @{ Row = 41; Price = '12.62'; Volume = '  -1.66%  ' }
    @{ Row = 42; Price = '0.6860'; Volume = '  -0.69%  ' }
    @{ Row = 43; Price = '1.270'; Volume = '  +2.11%  ' }
    @{ Row = 44; Price = '14.85'; Volume = '  +5.59%  ' }
    @{ Row = 45; Price = '0.6418'; Volume = '  +0.49%  ' }
    @{ Row = 46; Price = '2.312'; Volume = '  -1.82%  ' }
    @{ Row = 47; Price = '3.752'; Volume = '  +2.78%  ' }
    @{ Row = 48; Price = '0.00000000361'; Volume = '  +1.83%  ' }
    @{ Row = 49; Price = '1.254'; Volume = '  +0.45%  ' }
    @{ Row = 50; Price = '82.71'; Volume = '  -1.02%  ' }
    @{ Row = 51; Price = '0.07283'; Volume = '  +2.62%  ' }
)

foreach ($u in $updates) {
    $priceCell = $ws.Cells.Item($u.Row, 4)
    # Price column sometimes looks numeric ("54.12") or uses dotted
    # thousands grouping ("30.581.38") scraped verbatim as text --
    # force text storage so COM does not coerce it into a Double
    # and mangle the formatting (trailing zeros, sci notation, etc).
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $u.Price
    $priceCell.Style = "Normal"

    $ws.Cells.Item($u.Row, 5).Value = $u.Volume
}
